$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $cell = $ws.Range($cellRef)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue "D2" "36.503.31"
$ws.Range("E2").Value = "  -0.10%  "
Set-TextValue "D3" "1.949.39"
$ws.Range("E3").Value = "  +0.39%  "
$ws.Range("E4").Value = "  -0.05%  "
Set-TextValue "D5" "243.04"
$ws.Range("E5").Value = "  -0.29%  "
$ws.Range("E6").Value = "  +0.06%  "
Set-TextValue "D7" "60.41"
$ws.Range("E7").Value = "  +5.20%  "
$ws.Range("E8").Value = "  -0.01%  "
Set-TextValue "D9" "0.375"
$ws.Range("E9").Value = "  +3.49%  "
Set-TextValue "D10" "0.0785"
$ws.Range("E10").Value = "  -7.50%  "
$ws.Range("E11").Value = "  +0.37%  "
Set-TextValue "D12" "14.19"
$ws.Range("E12").Value = "  +5.64%  "
Set-TextValue "D13" "2.236.25"
$ws.Range("E13").Value = "  +0.38%  "
Set-TextValue "D14" "0.822"
$ws.Range("E14").Value = "  +1.22%  "
Set-TextValue "D15" "21.44"
$ws.Range("E15").Value = "  +0.43%  "
Set-TextValue "D16" "5.23"
$ws.Range("E16").Value = "  +1.14%  "
Set-TextValue "D17" "1.950.55"
$ws.Range("E17").Value = "  +0.51%  "
Set-TextValue "D18" "36.357.04"
$ws.Range("E18").Value = "  -0.28%  "
Set-TextValue "D19" "69.25"
$ws.Range("E19").Value = "  -0.22%  "
Set-TextValue "D20" "0.0₃0847"
$ws.Range("E20").Value = "  -2.06%  "
Set-TextValue "D21" "228.52"
$ws.Range("E21").Value = "  +0.10%  "
Set-TextValue "D22" "5.05"
$ws.Range("E22").Value = "  +1.11%  "
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("E24").Value = "  +3.29%  "
Set-TextValue "D25" "2.35"
$ws.Range("E25").Value = "  +2.03%  "
$ws.Range("E26").Value = "  +7.54%  "
Set-TextValue "D27" "9.14"
$ws.Range("E27").Value = "  -0.73%  "
Set-TextValue "D28" "159.70"
$ws.Range("E28").Value = "  -0.80%  "
Set-TextValue "D29" "19.23"
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("E30").Value = "  +18.74%  "
$ws.Range("E31").Value = "  +0.69%  "
$ws.Range("E32").Value = "  +2.78%  "
Set-TextValue "D33" "0.0610"
$ws.Range("E33").Value = "  -1.08%  "
Set-TextValue "D34" "4.42"
$ws.Range("E34").Value = "  +5.65%  "
$ws.Range("B35").Value = "RenderToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D35" "3.47"
$ws.Range("E35").Value = "  +9.03%  "
$ws.Range("B36").Value = "BinanceUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue "D36" "1.00"
$ws.Range("E36").Value = "  -0.20%  "
Set-TextValue "D37" "2.26"
$ws.Range("E37").Value = "  +3.93%  "
$ws.Range("E38").Value = "  -1.55%  "
$ws.Range("E39").Value = "  -13.19%  "
$ws.Range("E40").Value = "  +0.64%  "
Set-TextValue "D41" "0.0956"
$ws.Range("E41").Value = "  -2.90%  "
Set-TextValue "D42" "1.17"
$ws.Range("E42").Value = "  +1.35%  "
$ws.Range("E43").Value = "  -0.23%  "
Set-TextValue "D44" "1.355.79"
$ws.Range("E44").Value = "  +1.09%  "
$ws.Range("E45").Value = "  -2.12%  "
Set-TextValue "D46" "88.26"
$ws.Range("E46").Value = "  +1.82%  "
$ws.Range("E47").Value = "  -0.74%  "
Set-TextValue "D48" "7.09"
$ws.Range("E48").Value = "  -1.51%  "
$ws.Range("E49").Value = "  +0.52%  "
Set-TextValue "D50" "45.16"
$ws.Range("E50").Value = "  +4.45%  "
Set-TextValue "D51" "2.129.87"
$ws.Range("E51").Value = "  +0.50%  "
